$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-08-18 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-08-19 Monday", 2)

$d.Content.Find.Execute("91×72=", $true, $false, $false, $false, $false, $true, 1, $false, "61×24=", 2)
$d.Content.Find.Execute("14×47=", $true, $false, $false, $false, $false, $true, 1, $false, "76×47=", 2)
$d.Content.Find.Execute("61×64=", $true, $false, $false, $false, $false, $true, 1, $false, "79×75=", 2)
$d.Content.Find.Execute("51×28=", $true, $false, $false, $false, $false, $true, 1, $false, "32×26=", 2)
$d.Content.Find.Execute("22×98=", $true, $false, $false, $false, $false, $true, 1, $false, "18×14=", 2)

$d.Content.Find.Execute("16×28=", $true, $false, $false, $false, $false, $true, 1, $false, "64×22=", 2)
$d.Content.Find.Execute("51×22=", $true, $false, $false, $false, $false, $true, 1, $false, "99×22=", 2)
$d.Content.Find.Execute("18×47=", $true, $false, $false, $false, $false, $true, 1, $false, "55×65=", 2)
$d.Content.Find.Execute("61×95=", $true, $false, $false, $false, $false, $true, 1, $false, "90×39=", 2)
$d.Content.Find.Execute("21×31=", $true, $false, $false, $false, $false, $true, 1, $false, "92×66=", 2)

$d.Content.Find.Execute("65×27=", $true, $false, $false, $false, $false, $true, 1, $false, "19×14=", 2)
$d.Content.Find.Execute("94×70=", $true, $false, $false, $false, $false, $true, 1, $false, "67×44=", 2)
$d.Content.Find.Execute("18×39=", $true, $false, $false, $false, $false, $true, 1, $false, "50×11=", 2)
$d.Content.Find.Execute("31×37=", $true, $false, $false, $false, $false, $true, 1, $false, "47×61=", 2)
$d.Content.Find.Execute("12×80=", $true, $false, $false, $false, $false, $true, 1, $false, "69×54=", 2)

$d.Content.Find.Execute("58×99=", $true, $false, $false, $false, $false, $true, 1, $false, "34×14=", 2)
$d.Content.Find.Execute("95×31=", $true, $false, $false, $false, $false, $true, 1, $false, "94×88=", 2)
$d.Content.Find.Execute("99×31=", $true, $false, $false, $false, $false, $true, 1, $false, "39×60=", 2)
$d.Content.Find.Execute("57×35=", $true, $false, $false, $false, $false, $true, 1, $false, "57×33=", 2)
$d.Content.Find.Execute("80×52=", $true, $false, $false, $false, $false, $true, 1, $false, "27×94=", 2)

$d.Content.Find.Execute("97×43=", $true, $false, $false, $false, $false, $true, 1, $false, "41×48=", 2)
$d.Content.Find.Execute("72×59=", $true, $false, $false, $false, $false, $true, 1, $false, "99×31=", 2)
$d.Content.Find.Execute("28×34=", $true, $false, $false, $false, $false, $true, 1, $false, "63×33=", 2)
$d.Content.Find.Execute("30×55=", $true, $false, $false, $false, $false, $true, 1, $false, "69×70=", 2)
$d.Content.Find.Execute("99×47=", $true, $false, $false, $false, $false, $true, 1, $false, "98×27=", 2)
